$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.484.91'
$ws.Range("E2").Value = '  +0.79%  '

$ws.Range("D3").Value = '2.016.95'
$ws.Range("E3").Value = '  +0.70%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '262.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.619'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.82%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.18%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.385'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.16%  '

$ws.Range("E10").Value = '  -4.18%  '

$ws.Range("E11").Value = '  -1.63%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.34'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.05%  '

$ws.Range("D13").Value = '2.313.61'
$ws.Range("E13").Value = '  +0.71%  '

$ws.Range("E14").Value = '  -5.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.96'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -8.26%  '

$ws.Range("E16").Value = '  -4.10%  '

$ws.Range("D17").Value = '2.018.52'
$ws.Range("E17").Value = '  +0.31%  '

$ws.Range("D18").Value = '37.345.03'
$ws.Range("E18").Value = '  +0.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.68'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.09%  '

$ws.Range("E20").Value = '  -2.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.18'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.69'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.65%  '

$ws.Range("E24").Value = '  -0.06%  '

$ws.Range("E25").Value = '  -0.71%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.80'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.127'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -11.83%  '

$ws.Range("E30").Value = '  -1.42%  '

$ws.Range("E31").Value = '  -1.38%  '

$ws.Range("E32").Value = '  -3.58%  '

$ws.Range("E33").Value = '  -1.77%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.52'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.29%  '

$ws.Range("E35").Value = '  -0.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.83'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.79%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.17'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.03%  '

$ws.Range("E40").Value = '  +4.64%  '

$ws.Range("E41").Value = '  +3.87%  '

$ws.Range("E42").Value = '  -4.22%  '

$ws.Range("E43").Value = '  -0.72%  '

$ws.Range("D44").Value = '1.399.74'
$ws.Range("E44").Value = '  +1.37%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.28%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.72'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.65%  '

$ws.Range("E47").Value = '  -2.18%  '

$ws.Range("E48").Value = '  -4.41%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.90'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.84%  '

$ws.Range("D50").Value = '2.205.16'
$ws.Range("E50").Value = '  +0.71%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.96'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.26%  '

# Row 37/38 content swap (RenderToken <-> BinanceUSD reorder)
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.38'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.27%  '

$ws.Range("B38").Value = 'BinanceUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.10%  '
